# Generate Report for Handoff
# Updates the localization-status report with a new handoff id
# (374fccbe-4d42-4f2d-8204-80ce55b51c13 -> 0c20266b-14bd-43bc-ad79-9cc603363980)
# and refreshed handoff timestamps, across the Overview/zh-cn/de-de sheets.
# The hyperlink targets themselves (pointing at the historical commit blobs)
# are untouched - only the displayed file names / timestamps change.
#
# Note: this runtime's Range.Hyperlinks.Delete() clears every hyperlink on
# the worksheet (not just the ones touching that range), so each sheet's
# hyperlinks are dropped once and rebuilt in their original order.

$wb = $excel.ActiveWorkbook

$newId = "0c20266b-14bd-43bc-ad79-9cc603363980"
$newHash = "5f75b69751889da01ce1c6635672cab63374b8f5"

$newMdName = "$newId.md"
$newZhName = "$newId.$newHash.zh-cn.xlf"
$newDeName = "$newId.$newHash.de-de.xlf"

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/812f534927a7bb7573a4fcd90de5c3a92c2dc496/e2e/374fccbe-4d42-4f2d-8204-80ce55b51c13.md"
$zhTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67aceb04f011e4105ff3ed28bedc0f6dc31c03f7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/374fccbe-4d42-4f2d-8204-80ce55b51c13.02dd8d138a5f80cad5b4d516d71755f479ad6f8b.zh-cn.xlf"
$deTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6116b8a2e88ca72b3f2d775de9884f2a4bb24f10/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/374fccbe-4d42-4f2d-8204-80ce55b51c13.02dd8d138a5f80cad5b4d516d71755f479ad6f8b.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.UsedRange.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdTarget, "", "", $newMdName)
$wsOverview.Range("D2").Value = "2016-03-22 17:07:58"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.UsedRange.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdTarget, "", "", $newMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhTarget, "", "", $newZhName)
$wsZhCn.Range("E2").Value = "2016-03-22 17:07:55"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.UsedRange.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdTarget, "", "", $newMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deTarget, "", "", $newDeName)
$wsDeDe.Range("E2").Value = "2016-03-22 17:07:58"
